$d = $word.ActiveDocument

# Locate the anchor paragraph ("Requisitos" prerequisite line) and delete the
# three paragraphs that follow it: the blank spacer paragraph, the blank
# page-break paragraph, and the "© 2020 ..." footer paragraph - leaving the
# trailing blank paragraph + page-break paragraph untouched.
$anchorText = "LOQ4212: Engenharia da Qualidade II (Requisito fraco)"

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {
        $anchor = $p
        break
    }
}

$startPara = $anchor.Next()
$endPara = $startPara.Next().Next()

$start = $startPara.Range.Start
$end = $endPara.Range.End

$d.Range($start, $end).Delete()
